$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.256.40'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.596.87'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.244'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0854'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = '1.822.59'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '1.590.56'
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").Value = '  -2.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '26.257.94'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.37%  '
$ws.Range("E19").Value = '  +4.49%  '
$ws.Range("D20").Value = '0.0₃0721'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = '1.465.27'
$ws.Range("E33").Value = '  +3.44%  '
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.74%  '
$ws.Range("E37").Value = '  -3.12%  '
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  +2.27%  '
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("D44").Value = '1.735.12'
$ws.Range("E44").Value = '  +0.83%  '
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0501'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
